$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    $shape.TextFrame.WordWrap = $false
}
